$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17
$ws.Cells.Item($row, 1).Value = 42620.885636574072
$ws.Cells.Item($row, 2).Value = 22
$ws.Cells.Item($row, 3).Value = 58
$ws.Cells.Item($row, 4).Value = 38
$ws.Cells.Item($row, 5).Value = 58
$ws.Cells.Item($row, 6).Value = 6
$ws.Cells.Item($row, 7).Value = 26587
$ws.Cells.Item($row, 8).Value = 21209
$ws.Cells.Item($row, 9).Value = 1239
$ws.Cells.Item($row, 10).Value = 253
$ws.Cells.Item($row, 11).Value = 167
$ws.Cells.Item($row, 12).Value = 27
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = "Named"

$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
